$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# ---------------------------------------------------------------------------
# 1. Remove the empty "Sheet1" worksheet (4th tab) that was added by mistake.
# ---------------------------------------------------------------------------
$sheet1 = $wb.Worksheets.Item("Sheet1")
$sheet1.Delete()

# ---------------------------------------------------------------------------
# 2. Update the "ProII" worksheet (3rd tab) with the refreshed OID listing.
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("ProII")

# Widen column A to fit the longer class names (closest width the engine can
# persist to 37.54296875 given its internal pixel-grid quantization).
$ws.Columns.Item(1).ColumnWidth = 36.67

# Row 2
$ws.Range("A2").Value = "Simple HX"
$ws.Range("B2").Value = "HeatExchanger"
$ws.Range("C2").Value = 14

# Row 3
$ws.Range("A3").Value = "Pump"
$ws.Range("B3").Value = "Pump"
$ws.Range("C3").Value = 110

# Row 4
$ws.Range("A4").Value = "Flash"
$ws.Range("B4").Value = "Separator"
$ws.Range("C4").Value = 152

# Row 5
$ws.Range("A5").Value = "Valve"
$ws.Range("B5").Value = "Valve"
$ws.Range("C5").Value = 222

# Row 6
$ws.Range("A6").Value = "Compressor"
$ws.Range("B6").Value = "Compressor"
$ws.Range("C6").Value = 269

# Row 7
$ws.Range("A7").Value = "Air Cooled HX"
$ws.Range("B7").Value = "AirCooledExchanger"
$ws.Range("C7").Value = 317

# Row 8
$ws.Range("A8").Value = "Stream"
$ws.Range("B8").Value = "PipingSystem"
$ws.Range("C8").Value = 412

# Row 9
$ws.Range("A9").Value = "Bulk Phase"
$ws.Range("B9").Value = "PipingSystem"
$ws.Range("C9").Value = 430

# Row 10
$ws.Range("A10").Value = "Liquid Phase"
$ws.Range("B10").Value = "PipingSystem"
$ws.Range("C10").Value = 430

# Row 11
$ws.Range("A11").Value = "Vapor Phase"
$ws.Range("B11").Value = "PipingSystem"
$ws.Range("C11").Value = 448

# Row 14
$ws.Range("A14").Value = "Distilation Column and Trays"
$ws.Range("B14").Value = "Distillation"
$ws.Range("C14").Value = 596

# Rows 15-19 (new rows) are filled in before rows 12/13/20 are revisited below,
# matching the order in which the author actually introduced the new class
# names (this controls the append order of the shared-string table).

# Row 15 (new)
$ws.Range("A15").Value = "PumpAround"
$ws.Range("B15").Value = "Distillation"
$ws.Range("C15").Value = 729
$ws.Range("C15").Font.Name = "Arial"

# Row 16 (new)
$ws.Range("A16").Value = "Tubine"
$ws.Range("B16").Value = "Expander"
$ws.Range("C16").Value = 949
$ws.Range("C16").Font.Name = "Arial"

# Row 17 (new)
$ws.Range("A17").Value = "SideColumn"
$ws.Range("B17").Value = "Distillation"
$ws.Range("C17").Value = 996
$ws.Range("C17").Font.Name = "Arial"

# Row 18 (new)
$ws.Range("A18").Value = "DistillationPortData"
$ws.Range("B18").Value = "Distillation"

# Row 19 (new)
$ws.Range("A19").Value = "ValvePortData"
$ws.Range("B19").Value = "Valve"

# Row 12
$ws.Range("A12").Value = "Mixer"
$ws.Range("B12").Value = "Generic"
$ws.Range("C12").Value = 474

# Row 13
$ws.Range("A13").Value = "Splitter"
$ws.Range("B13").Value = "Generic"
$ws.Range("C13").Value = 535

# Row 20 (new)
$ws.Range("A20").Value = "RigorousExchanger"
$ws.Range("B20").Value = "Expander"

# Restore the tab selection to the cell the author left active.
$ws.Range("C18").Select()
